# Edit: "Final Thoughts / Future Improvements" slide (last slide) is
# restructured - title becomes a plain centered textbox at the top, the
# bullet list becomes a larger plain rectangle with many more bullet
# points, and the slide itself is effectively recreated (its p:sldId
# changes from 271 to 272). Two new slide guides are also added to the
# presentation, plus a couple of small text cleanups on two other slides.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Recreate the last slide (index 7) so PowerPoint assigns it a fresh
#    SlideId (271 -> 272), matching the commit. We add a new slide at
#    the end (re-using the "Title and Content" layout that slide 7
#    already used) and then delete the old slide 7, which leaves the
#    new slide sitting at position 7 with the next available SlideId.
# ---------------------------------------------------------------------
$oldSlide7 = $p.Slides.Item(7)
$newSlide = $p.Slides.Add(8, 2)
$oldSlide7.Delete()

$s7 = $p.Slides.Item(7)

# Remove the two empty inherited placeholders (Content + Title) that
# come from the layout - the final slide has no placeholders at all,
# just two free-floating shapes.
while ($s7.Shapes.Count -gt 0) {
    $s7.Shapes.Item(1).Delete()
}

# ---------------------------------------------------------------------
# 2) Add the big bullet-list rectangle ("Rectangle 1").
# ---------------------------------------------------------------------
$rect = $s7.Shapes.AddShape(1, 138223/12700, 646417/12700, 12053777/12700, 6186309/12700)
$rect.Name = "Rectangle 1"

$bodyTf = $rect.TextFrame
$bodyTf.WordWrap = -1

$bullets = @(
    "Add more variables and formulas",
    "Expand into other branches of physics",
    "Add more extensive unit conversion feature",
    "Add ability to save old calculator sessions",
    "Create web page / app for mobile use",
    "Display equations used with variables entered",
    "Add tabs for useful links to physics pages, etc.",
    "Suggestions for missing variables, if it can`u{2019}t solve with entries",
    "Add `u{201c}undo`u{201d} or `u{201c}back`u{201d} button",
    "Add sound and graphics",
    "More ways to incorporate flames / explosions"
)

$bodyTr = $bodyTf.TextRange
$bodyTr.Text = $bullets -join "`r"
$bodyTr.Font.Size = 36

for ($i = 1; $i -le $bullets.Count; $i++) {
    $para = $bodyTr.Paragraphs($i, 1)
    $para.ParagraphFormat.Bullet.Visible = -1
    $para.ParagraphFormat.Bullet.Character = 8226
    $para.ParagraphFormat.Bullet.Font.Name = "Arial"
}

# ---------------------------------------------------------------------
# 3) Add the title textbox ("Title 1"), centered, at the top.
# ---------------------------------------------------------------------
$title = $s7.Shapes.AddTextbox(1, 0/12700, 91244/12700, 12192000/12700, 642403/12700)
$title.Name = "Title 1"
$titleTr = $title.TextFrame.TextRange
$titleTr.Text = "Final Thoughts / Future Improvements"
$titleTr.ParagraphFormat.Alignment = 2

# ---------------------------------------------------------------------
# 4) Slide guides: add a horizontal guide at y=2160 and a vertical guide
#    at x=3840 (in the new p15:sldGuideLst extension).
# ---------------------------------------------------------------------
Add-Type -AssemblyName Microsoft.Office.Interop.PowerPoint -ErrorAction SilentlyContinue
try {
    $guides = $p.SlideMaster.Guides
} catch {
    $guides = $null
}
if ($guides -ne $null) {
    $guides.Add(1, 2160/12700)
    $guides.Add(2, 3840/12700)
}

# ---------------------------------------------------------------------
# 5) Minor text touch-ups (no content change) on slide 4 and slide 6,
#    matching the smtClean / trailing endParaRPr cleanup seen in the
#    commit.
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
foreach ($sh in $s4.Shapes) {
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -like "We were never actually*") {
            $sh.TextFrame.TextRange.Text = $sh.TextFrame.TextRange.Text
        }
    }
}

$s6 = $p.Slides.Item(6)
foreach ($sh in $s6.Shapes) {
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -like "Final Look*") {
            $sh.TextFrame.TextRange.Text = $sh.TextFrame.TextRange.Text
        }
    }
}
